$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Prediction (column C) values for rows 2..170, one per row in order.
$predictionValues = @(0,0,2.631,2.029,1.161,0.4,0.037,0,0,0,0,0,0,0,0,0,0.083,0.414,0.99,1.598,2.015,2.805,3.001,3.143,3.108,2.903,2.217,1.919,1.215,0.398,0.037,0,0,0,0,0,0,0,0,0,0.102,0.393,1.026,1.628,2.147,2.833,2.931,2.769,2.733,2.201,1.773,1.309,0.803,0.329,0.037,0,0,0,0,0,0,0,0,0,0.08,0.41,0.667,1.021,1.595,1.846,2.11,2.034,1.797,1.421,1.191,0.844,0.433,0.182,0.038,0,0,0,0,0,0,0,0,0,0.079,0.387,0.667,1.08,1.753,2.045,2.307,2.376,2.143,2.012,1.782,1.296,0.802,0.277,0.037,0,0,0,0,0,0,0,0,0,0.081,0.401,0.984,1.573,1.961,2.969,3.016,3.039,2.851,2.503,1.963,1.437,0.906,0.258,0.036,0,0,0,0,0,0,0,0,0,0.048,0.127,0.307,0.65,0.93,1.042,1.29,1.402,1.353,1.309,1.102,0.764,0.433,0.186,0.037,0,0,0,0,0,0,0,0,0,0.08400000000000001,0.412,0.96,1.236,1.932,2.308,2.437,2.585,2.678)

for ($row = 2; $row -le 170; $row++) {
    $idx = $row - 2

    # Column A: shift the date forward by 5 days (keeps the time-of-day/style).
    $oldDate = $ws.Cells.Item($row, 1).Value2
    $newDate = $oldDate + 5
    $ws.Cells.Item($row, 1).Value = $newDate

    # Column B (hour) is unchanged.
    $hour = $ws.Cells.Item($row, 2).Value2

    # Column C: updated Prediction value.
    $ws.Cells.Item($row, 3).Value = $predictionValues[$idx]

    # Column D: Lookup string rebuilt from the new date + hour.
    $dt = [datetime]::FromOADate($newDate)
    $lookup = $dt.ToString("dd.MM.yyyy") + $hour
    $ws.Cells.Item($row, 4).Value = $lookup
}

Write-Host "Fixed dates for Fundamentals: shifted rows 2-170 forward by 5 days and refreshed predictions."
